$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B64").Value = "Film de fusor amarillo HP"
$ws.Range("D64").Value = 10000
$ws.Range("E64").Value = 70000
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 12
$ws.Range("J64").Value = 0
